$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3454506666666666
$ws.Range("H2").Value = 1.036352
$ws.Range("I2").Value = 0.1052716477644991
$ws.Range("J2").Value = 0.1052716477644991
$ws.Range("M2").Value = 15.67027366666667
$ws.Range("N2").Value = 47.010821
$ws.Range("O2").Value = 0.3497296100206518
$ws.Range("P2").Value = 0.3497296100206518
$ws.Range("Q2").Value = 5.413306484999111
$ws.Range("R2").Value = 48.719758364992
$ws.Range("S2").Value = 0.0368166123189097
$ws.Range("T2").Value = 0.03681661231890971
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3454506666666666
$ws.Range("H3").Value = 1.036352
$ws.Range("I3").Value = 0.1052716477644991
$ws.Range("J3").Value = 0.1052716477644991
$ws.Range("O3").Value = 0.5842073451353695
$ws.Range("P3").Value = 0.5842073451353695
$ws.Range("Q3").Value = 9.042681315484444
$ws.Range("R3").Value = 81.38413183935999
$ws.Range("S3").Value = 0.0615004698585238
$ws.Range("T3").Value = 0.0615004698585238
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3454506666666666
$ws.Range("H4").Value = 1.036352
$ws.Range("I4").Value = 0.1052716477644991
$ws.Range("J4").Value = 0.1052716477644991
$ws.Range("O4").Value = 0.06606304484397865
$ws.Range("P4").Value = 0.06606304484397864
$ws.Range("Q4").Value = 1.022559997283556
$ws.Range("R4").Value = 9.203039975552
$ws.Range("S4").Value = 0.006954565587065631
$ws.Range("T4").Value = 0.006954565587065631
$ws.Range("I5").Value = 0.3398937483175971
$ws.Range("J5").Value = 0.3398937483175971
$ws.Range("M5").Value = 15.67027366666667
$ws.Range("N5").Value = 47.010821
$ws.Range("O5").Value = 0.3497296100206518
$ws.Range("P5").Value = 0.3497296100206518
$ws.Range("Q5").Value = 17.478106128769
$ws.Range("R5").Value = 157.302955158921
$ws.Range("S5").Value = 0.1188709080475708
$ws.Range("T5").Value = 0.1188709080475708
$ws.Range("I6").Value = 0.3398937483175971
$ws.Range("J6").Value = 0.3398937483175971
$ws.Range("O6").Value = 0.5842073451353695
$ws.Range("P6").Value = 0.5842073451353695
$ws.Range("S6").Value = 0.1985684243327329
$ws.Range("T6").Value = 0.1985684243327329
$ws.Range("I7").Value = 0.3398937483175971
$ws.Range("J7").Value = 0.3398937483175971
$ws.Range("O7").Value = 0.06606304484397865
$ws.Range("P7").Value = 0.06606304484397864
$ws.Range("S7").Value = 0.02245441593729341
$ws.Range("T7").Value = 0.0224544159372934
$ws.Range("I8").Value = 0.5548346039179038
$ws.Range("J8").Value = 0.5548346039179038
$ws.Range("M8").Value = 15.67027366666667
$ws.Range("N8").Value = 47.010821
$ws.Range("O8").Value = 0.3497296100206518
$ws.Range("P8").Value = 0.3497296100206518
$ws.Range("Q8").Value = 28.53085159462633
$ws.Range("R8").Value = 256.777664351637
$ws.Range("S8").Value = 0.1940420896541713
$ws.Range("T8").Value = 0.1940420896541713
$ws.Range("I9").Value = 0.5548346039179038
$ws.Range("J9").Value = 0.5548346039179038
$ws.Range("O9").Value = 0.5842073451353695
$ws.Range("P9").Value = 0.5842073451353695
$ws.Range("S9").Value = 0.3241384509441129
$ws.Range("T9").Value = 0.3241384509441129
$ws.Range("I10").Value = 0.5548346039179038
$ws.Range("J10").Value = 0.5548346039179038
$ws.Range("O10").Value = 0.06606304484397865
$ws.Range("P10").Value = 0.06606304484397864
$ws.Range("S10").Value = 0.03665406331961962
$ws.Range("T10").Value = 0.03665406331961961
